$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove the placeholder "Sheet2" (X1/X2/X3) now that real participants
# are tracked on Sheet1 -------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Delete()

$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Correct Aniket Saha's email address --------------------------------
$ws1.Range("B2").Value = "anikets349@gmail.com"

# --- Add the two new participants for the email workflow ---------------
$ws1.Range("A6").Value = "Madhumita Saha"
$ws1.Range("B6").Value = "madhumitasaha898@gmail.com"
$ws1.Range("A7").Value = "Abhijit Saha"
$ws1.Range("B7").Value = "abhijitsaha8698@gmail.com"

# Drop the stale hyperlinks first - they'll be rebuilt after the sort so
# each one tracks the row its own address ends up on.
$ws1.Cells.Hyperlinks.Delete()

# --- Sort the participant table A-Z by name (header stays put) ---------
$sortRange = $ws1.Range("A2:B7")
$sortRange.Sort($ws1.Range("A2:A7"))

# --- Re-create the mailto hyperlinks against the sorted rows ------------
# (read each row's own email back so the link always matches the address
# that landed there, regardless of how the sort reordered things)
for ($r = 2; $r -le 7; $r++) {
    $email = $ws1.Cells.Item($r, 2).Value2
    $cell = $ws1.Cells.Item($r, 2)
    $ws1.Hyperlinks.Add($cell, "mailto:" + $email)
}

# --- Cosmetic tweaks that came along with the edit ----------------------
$ws1.Columns.Item(1).ColumnWidth = 17.16
$ws1.Columns.Item(2).ColumnWidth = 27.16

$ws1.Application.ActiveWindow.Zoom = 110

$ws1.Range("F8").Select()
